$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'55.260.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.36%  "
# Row 3
$ws.Range("D3").Value = "'2.352.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.56%  "
# Row 4
$ws.Range("E4").Value = "  -0.02%  "
# Row 5
$ws.Range("D5").Value = "'476.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.21%  "
# Row 6
$ws.Range("D6").Value = "'146.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.16%  "
# Row 7
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "
# Row 8
$ws.Range("E8").Value = "  +20.58%  "
# Row 9
$ws.Range("D9").Value = "'2.356.21"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.60%  "
# Row 10
$ws.Range("D10").Value = "'0.0961"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.66%  "
# Row 11
$ws.Range("D11").Value = "'5.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.04%  "
# Row 12
$ws.Range("E12").Value = "  -1.62%  "
# Row 13
$ws.Range("E13").Value = "  +1.37%  "
# Row 14
$ws.Range("D14").Value = "'2.762.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.53%  "
# Row 15
$ws.Range("D15").Value = "'55.224.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.47%  "
# Row 16
$ws.Range("D16").Value = "'20.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.84%  "
# Row 17
$ws.Range("E17").Value = "  -4.52%  "
# Row 18
$ws.Range("D18").Value = "'2.353.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.71%  "
# Row 19
$ws.Range("D19").Value = "'4.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.31%  "
# Row 20
$ws.Range("D20").Value = "'314.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.59%  "
# Row 21
$ws.Range("D21").Value = "'9.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.40%  "
# Row 22
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.23%  "
# Row 23
$ws.Range("D23").Value = "'5.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.84%  "
# Row 24
$ws.Range("D24").Value = "'56.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.47%  "
# Row 25
$ws.Range("E25").Value = "  +0.03%  "
# Row 26
$ws.Range("D26").Value = "'0.395"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.70%  "
# Row 27
$ws.Range("D27").Value = "'0.153"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.82%  "
# Row 28
$ws.Range("D28").Value = "'2.444.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.05%  "
# Row 29
$ws.Range("D29").Value = "'7.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.86%  "
# Row 30
$ws.Range("E30").Value = "  +0.09%  "
# Row 31
$ws.Range("D31").Value = "'0.0₃0743"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.75%  "
# Row 32
$ws.Range("D32").Value = "'145.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.44%  "
# Row 33
$ws.Range("D33").Value = "'18.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.70%  "
# Row 34
$ws.Range("E34").Value = "  -1.87%  "
# Row 35
$ws.Range("D35").Value = "'5.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.04%  "
# Row 36
$ws.Range("D36").Value = "'1.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.90%  "
# Row 37
$ws.Range("D37").Value = "'3.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.18%  "
# Row 38
$ws.Range("D38").Value = "'0.812"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.36%  "
# Row 39
$ws.Range("D39").Value = "'0.102"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.66%  "
# Row 40
$ws.Range("D40").Value = "'33.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.33%  "
# Row 42
$ws.Range("E42").Value = "  -0.01%  "
# Row 43
$ws.Range("D43").Value = "'3.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.66%  "
# Row 44
$ws.Range("D44").Value = "'0.577"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.32%  "
# Row 45
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0516"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.90%  "
# Row 46
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "'10.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.44%  "
# Row 47
$ws.Range("D47").Value = "'249.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.38%  "
# Row 48
$ws.Range("D48").Value = "'0.0221"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.58%  "
# Row 49
$ws.Range("D49").Value = "'4.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.17%  "
# Row 50
$ws.Range("D50").Value = "'1.794.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.10%  "
# Row 51
$ws.Range("D51").Value = "'16.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.76%  "
